$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 113 (Item #112): fill in the "E" disposition with Complete and restyle
# the row to match the "filled" banding used by the rest of the completed
# items (copy formatting from row 108, which already uses that style group).
$src = $ws.Range("A108:E108")
$dst = $ws.Range("A113:E113")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B113").Value = "Review caIntegrator performance improvements in the 7/23 status meeting."
$ws.Range("C113").Value = "Mike Hunter and Abe Evans-El"
$ws.Range("D113").Value = 40009
$ws.Range("E113").Value = "Complete"

# Row 114 (Item #113): new action item
$ws.Range("B114").Value = "Review the project plans for caArray and caIntegrator to determine if there are things to do to streamline getting both applications back into Production before September"
$ws.Range("C114").Value = "JJ Pan"
$ws.Range("D114").Value = 40016
$ws.Range("E114").Value = "Assigned"
$ws.Rows.Item(114).RowHeight = 46

# Row 115 (Item #114): new action item
$ws.Range("B115").Value = "Decide which of the three options proposed for making it easy to review audit log entries each month to implement"
$ws.Range("C115").Value = "JJ Pan and Juli Klemm"
$ws.Range("D115").Value = 40016
$ws.Range("E115").Value = "Assigned"
$ws.Rows.Item(115).RowHeight = 31

# Row 116 (Item #115): new action item
$ws.Range("B116").Value = "After Abe addresses the performance drop for Agilent Data Sets, share Abe's performance page from the wiki with Eve Shalley"
$ws.Range("C116").Value = "Mike Hunter"
$ws.Range("D116").Value = 40016
$ws.Range("E116").Value = "Assigned"
$ws.Rows.Item(116).RowHeight = 31

# Update the window view: scroll so row 96 is at the top and select E117,
# matching where the status meeting left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 96
$win.ScrollColumn = 2
$ws.Range("E117").Select()
